$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.753.79'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.321.95'
$ws.Range("E3").Value = '  +3.77%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '96.88'
$ws.Range("E5").Value = '  +9.81%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '271.13'
$ws.Range("E6").Value = '  -0.54%  '
$ws.Range("E7").Value = '  +0.82%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.624'
$ws.Range("E9").Value = '  +3.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.45'
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("E11").Value = '  +2.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.11'
$ws.Range("E12").Value = '  +5.93%  '
$ws.Range("E13").Value = '  +0.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.668.51'
$ws.Range("E14").Value = '  +3.84%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.59'
$ws.Range("E15").Value = '  +3.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.857'
$ws.Range("E16").Value = '  +7.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.330.06'
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.685.96'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000109'
$ws.Range("E19").Value = '  +4.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.34'
$ws.Range("E20").Value = '  +6.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.85'
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '238.30'
$ws.Range("E22").Value = '  +2.48%  '
$ws.Range("E23").Value = '  -2.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.53'
$ws.Range("E24").Value = '  +9.72%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.56'
$ws.Range("E26").Value = '  +0.73%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.36'
$ws.Range("E27").Value = '  +4.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.43'
$ws.Range("E28").Value = '  -2.44%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.29'
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.86'
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("E31").Value = '  +8.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '172.03'
$ws.Range("E32").Value = '  -0.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0895'
$ws.Range("E33").Value = '  -2.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.50'
$ws.Range("E34").Value = '  +3.13%  '
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0358'
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.41'
$ws.Range("E37").Value = '  +3.31%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.107'
$ws.Range("E38").Value = '  -3.56%  '
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.35'
$ws.Range("E40").Value = '  +8.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.236'
$ws.Range("E41").Value = '  +10.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.36'
$ws.Range("E42").Value = '  +19.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.14'
$ws.Range("E43").Value = '  -2.29%  '
$ws.Range("E44").Value = '  -2.38%  '
$ws.Range("E45").Value = '  +7.68%  '
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("E47").Value = '  +4.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '100.56'
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("E49").Value = '  +2.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.549.40'
$ws.Range("E50").Value = '  +3.94%  '
$ws.Range("E51").Value = '  +13.07%  '
